$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric-looking strings (e.g. "5.44")
# must be forced to Text format first, so Excel keeps them as text
# (matching the original inline-string / "General" text formatting)
# instead of silently converting them into real numbers.
$textCells = @("D5","D6","D10","D12","D16","D17","D19","D21","D22","D24","D25","D27","D28","D29","D32","D34","D35","D36","D37","D39","D40","D41","D42","D43","D44","D46","D47","D48","D49","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "59.208.53"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "2.521.70"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "537.27"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").Value = "140.27"
$ws.Range("E6").Value = "  -2.13%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("E8").Value = "  -1.60%  "
$ws.Range("D9").Value = "2.525.15"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").Value = "0.100"
$ws.Range("E10").Value = "  +1.01%  "
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("D12").Value = "5.44"
$ws.Range("E12").Value = "  -1.75%  "
$ws.Range("E13").Value = "  +1.71%  "
$ws.Range("D14").Value = "2.966.95"
$ws.Range("E14").Value = "  +0.89%  "
$ws.Range("D15").Value = "59.144.93"
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("D16").Value = "22.92"
$ws.Range("E16").Value = "  -1.93%  "
$ws.Range("D17").Value = "0.0000141"
$ws.Range("E17").Value = "  +1.80%  "
$ws.Range("D18").Value = "2.543.58"
$ws.Range("E18").Value = "  +1.12%  "
$ws.Range("D19").Value = "10.92"
$ws.Range("E19").Value = "  -2.33%  "
$ws.Range("E20").Value = "  -0.62%  "
$ws.Range("D21").Value = "322.37"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("E23").Value = "  +1.63%  "
$ws.Range("D24").Value = "62.14"
$ws.Range("E24").Value = "  +0.63%  "
$ws.Range("D25").Value = "0.423"
$ws.Range("E26").Value = "  +1.48%  "
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").Value = "7.78"
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("D29").Value = "6.74"
$ws.Range("E29").Value = "  +1.16%  "
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").Value = "0.0₃0767"
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("D32").Value = "161.29"
$ws.Range("E32").Value = "  +1.85%  "
$ws.Range("E33").Value = "  +0.24%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "1.47"
$ws.Range("E34").Value = "  +3.43%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").Value = "1.13"
$ws.Range("E35").Value = "  -4.49%  "
$ws.Range("D36").Value = "18.53"
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("D37").Value = "4.21"
$ws.Range("E37").Value = "  -2.64%  "
$ws.Range("E38").Value = "  -1.17%  "
$ws.Range("D39").Value = "36.97"
$ws.Range("E39").Value = "  +1.53%  "
$ws.Range("D40").Value = "3.64"
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("D41").Value = "0.806"
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "284.40"
$ws.Range("E42").Value = "  -3.98%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "5.25"
$ws.Range("E43").Value = "  -5.04%  "
$ws.Range("D44").Value = "0.997"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("E45").Value = "  +1.09%  "
$ws.Range("D46").Value = "0.596"
$ws.Range("E46").Value = "  -1.29%  "
$ws.Range("D47").Value = "0.0931"
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").Value = "122.46"
$ws.Range("E48").Value = "  -1.64%  "
$ws.Range("D49").Value = "18.54"
$ws.Range("E49").Value = "  -0.40%  "
$ws.Range("D50").Value = "0.0511"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("E51").Value = "  -1.30%  "
